$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so values like "1.003" are not
# auto-converted to numbers by Excel's type inference.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.242.52"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "1.797.36"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "339.37"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.4576"
$ws.Range("E7").Value = "  +21.59%  "
$ws.Range("D8").Value = "0.3610"
$ws.Range("E8").Value = "  +7.32%  "
$ws.Range("D9").Value = "45.51"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").Value = "1.141"
$ws.Range("E10").Value = "  +1.24%  "
$ws.Range("D11").Value = "0.07529"
$ws.Range("E11").Value = "  +4.92%  "
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "22.40"
$ws.Range("E13").Value = "  +0.10%  "
$ws.Range("D14").Value = "6.222"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "7.250"
$ws.Range("E15").Value = "  +1.04%  "
$ws.Range("D16").Value = "1.792.57"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "0.00001084"
$ws.Range("E17").Value = "  +3.15%  "
$ws.Range("D18").Value = "0.06708"
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("D19").Value = "81.29"
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("D20").Value = "0.9992"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "17.23"
$ws.Range("E21").Value = "  +2.04%  "
$ws.Range("D22").Value = "6.374"
$ws.Range("E22").Value = "  +1.79%  "
$ws.Range("D23").Value = "28.220.90"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("D24").Value = "11.89"
$ws.Range("E24").Value = "  +1.94%  "
$ws.Range("D25").Value = "2.383"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("D26").Value = "20.38"
$ws.Range("E26").Value = "  +3.38%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.393"
$ws.Range("E27").Value = "  +3.11%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "153.45"
$ws.Range("E28").Value = "  +0.54%  "
$ws.Range("D29").Value = "1.997.53"
$ws.Range("E29").Value = "  +1.76%  "
$ws.Range("D30").Value = "1.270"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").Value = "132.34"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("D32").Value = "4.075"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("D33").Value = "5.880"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").Value = "0.09464"
$ws.Range("E34").Value = "  +7.89%  "
$ws.Range("D35").Value = "0.02377"
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("D36").Value = "12.07"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").Value = "0.06270"
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "0.6617"
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("D39").Value = "5.177"
$ws.Range("E39").Value = "  +1.20%  "
$ws.Range("D40").Value = "0.2158"
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("D41").Value = "1.481"
$ws.Range("E41").Value = "  +1.52%  "
$ws.Range("D42").Value = "1.213"
$ws.Range("E42").Value = "  +0.47%  "
$ws.Range("D43").Value = "8.064"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Value = "0.9990"
$ws.Range("E44").Value = "  -0.10%  "
$ws.Range("D45").Value = "13.95"
$ws.Range("E45").Value = "  +2.32%  "
$ws.Range("D46").Value = "3.873"
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").Value = "0.6074"
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("D48").Value = "128.34"
$ws.Range("E48").Value = "  -0.89%  "
$ws.Range("D49").Value = "2.025"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("D50").Value = "0.07094"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").Value = "1.166"
$ws.Range("E51").Value = "  -0.78%  "
